$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("Q2").Value = 3.4
$ws.Range("R2").Value = 1.33

# Row 3 updates
$ws.Range("G3").Value = 3
$ws.Range("I3").Value = 2.32
$ws.Range("J3").Value = 3.45
$ws.Range("L3").Value = 2.87
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.75
$ws.Range("W3").Value = 9.25
$ws.Range("X3").Value = 16
$ws.Range("AG3").Value = 7.8
$ws.Range("AH3").Value = 11.5
$ws.Range("AJ3").Value = 25
$ws.Range("AQ3").Value = 70
$ws.Range("AV3").Value = 4.3
$ws.Range("AW3").Value = 12
$ws.Range("AY3").Value = 50
